# "Add files via upload" - append 12 new login/register test rows (A15:B26)
# to the "login" sheet, and nudge the two sheets' remembered cell selections
# (login -> A2, register -> G9) the way re-saving the workbook in Excel did.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("login")
$ws2 = $wb.Worksheets.Item("register")

# New rows are written in the order that first introduces each new distinct
# string value (ad min / 326831XUliang / ad ddd / admin / aa aaa /
# 326831 XUliang) so the workbook's shared-string table grows the same way
# it did in the source edit (row 20 was keyed in just ahead of row 19).

$ws1.Range("A15").Value = "ad min"
$ws1.Range("B15").Value = "326831XUliang"

$ws1.Range("A16").Value = "ad min"
$ws1.Range("B16").Value = "aaaaa"

$ws1.Range("A17").Value = "ad ddd"
$ws1.Range("B17").Value = "326831XUliang"

$ws1.Range("A18").Value = "ad ddd"
$ws1.Range("B18").Value = "aaaaa"

$ws1.Range("A20").Value = "admin"
$ws1.Range("B20").Value = "aa aaa"

$ws1.Range("A19").Value = "admin"
$ws1.Range("B19").Value = "326831 XUliang"

$ws1.Range("A21").Value = "ad min"
$ws1.Range("B21").Value = "326831 XUliang"

$ws1.Range("A22").Value = "ad min"
$ws1.Range("B22").Value = "aa aaa"

$ws1.Range("A23").Value = "adddd"
$ws1.Range("B23").Value = "326831 XUliang"

$ws1.Range("A24").Value = "adddd"
$ws1.Range("B24").Value = "aa aaa"

$ws1.Range("A25").Value = "ad ddd"
$ws1.Range("B25").Value = "326831 XUliang"

$ws1.Range("A26").Value = "ad ddd"
$ws1.Range("B26").Value = "aa aaa"

# Restore/adjust the remembered selection on each sheet without disturbing
# which tab is active (register stays the active tab, as in the source file).
$ws1.Range("A2").Select()

$ws2.Activate()
$ws2.Range("G9").Select()
